$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)
$ws.Range("A5").Value = "Today 5"
